# "Removed vk link from front slide"
#
# The front slide (slide 1) has a small textbox ("TextBox 7") that just
# contains the link "http://vk.com/club33848893". Remove that shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$vkShape = $s.Shapes.Item("TextBox 7")
$vkShape.Delete()
